$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4312.968
$ws.Range("J17").Value = 4356.7334
$ws.Range("L17").Value = 13070.2002
$ws.Range("N17").Value = -13406.2002

$ws.Range("H76").Value = 4998
$ws.Range("I76").Value = 4998
$ws.Range("K76").Value = 4998
$ws.Range("M76").Value = -4683

$ws.Range("H79").Value = 4998
$ws.Range("I79").Value = 4998
$ws.Range("K79").Value = 4998
$ws.Range("M79").Value = -3906

$ws.Range("H98").Value = 3001
$ws.Range("I98").Value = 901.8
$ws.Range("J98").Value = 6499.6665
$ws.Range("K98").Value = 901.8
$ws.Range("L98").Value = 6499.6665
$ws.Range("M98").Value = 596.2
$ws.Range("N98").Value = -9495.666499999999

$ws.Range("H103").Value = 1398.8
$ws.Range("I103").Value = 999
$ws.Range("J103").Value = 1498.75
$ws.Range("K103").Value = 2997
$ws.Range("L103").Value = 4496.25
$ws.Range("M103").Value = -2411
$ws.Range("N103").Value = -5668.25

$ws.Range("H113").Value = 2988
$ws.Range("I113").Value = 2791.2
$ws.Range("J113").Value = 3972
$ws.Range("K113").Value = 2791.2
$ws.Range("L113").Value = 3972
$ws.Range("M113").Value = 462.8000000000002
$ws.Range("N113").Value = -10480

$ws.Range("H116").Value = 7586.5
$ws.Range("I116").Value = 6878.6
$ws.Range("J116").Value = 7908.273
$ws.Range("K116").Value = 6878.6
$ws.Range("L116").Value = 7908.273
$ws.Range("M116").Value = -3436.6
$ws.Range("N116").Value = -14792.273

$ws.Range("H122").Value = 3001
$ws.Range("I122").Value = 901.8
$ws.Range("J122").Value = 6499.6665
$ws.Range("K122").Value = 2705.4
$ws.Range("L122").Value = 19498.9995
$ws.Range("M122").Value = -255.3999999999996
$ws.Range("N122").Value = -24398.9995

$ws.Range("H132").Value = 790.21875
$ws.Range("I132").Value = 743.1
$ws.Range("K132").Value = 2229.3
$ws.Range("M132").Value = 300.6999999999998

$ws.Range("H135").Value = 588.4545000000001
$ws.Range("I135").Value = 588.4545000000001
$ws.Range("K135").Value = 5296.0905
$ws.Range("M135").Value = -2761.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 676666.7
$ws.Range("I8").Value = 676666.7
$ws.Range("K8").Value = 676666.7
$ws.Range("M8").Value = -676522.7

$ws.Range("H32").Value = 5363.5
$ws.Range("I32").Value = 5363.5
$ws.Range("K32").Value = 5363.5
$ws.Range("M32").Value = -5076.5

$ws.Range("H123").Value = 44950
$ws.Range("J123").Value = 44950
$ws.Range("L123").Value = 44950
$ws.Range("N123").Value = -54750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H94").Value = 4723.2
$ws.Range("I94").Value = 5313.923
$ws.Range("J94").Value = 883.5
$ws.Range("K94").Value = 5313.923
$ws.Range("L94").Value = 883.5
$ws.Range("M94").Value = -4862.923
$ws.Range("N94").Value = -1785.5

$ws.Range("H134").Value = 1630.95
$ws.Range("I134").Value = 1653.9474
$ws.Range("J134").Value = 1194
$ws.Range("K134").Value = 4961.8422
$ws.Range("L134").Value = 3582
$ws.Range("M134").Value = -2426.8422
$ws.Range("N134").Value = -8652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 10770.75
$ws.Range("J43").Value = 12694.333
$ws.Range("L43").Value = 12694.333
$ws.Range("N43").Value = -13062.333

$ws.Range("H99").Value = 3505.3333
$ws.Range("I99").Value = 1185
$ws.Range("K99").Value = 1185
$ws.Range("M99").Value = 313

$ws.Range("H101").Value = 10770.75
$ws.Range("J101").Value = 12694.333
$ws.Range("L101").Value = 12694.333
$ws.Range("N101").Value = -19184.333

$ws.Range("H126").Value = 3505.3333
$ws.Range("I126").Value = 1185
$ws.Range("K126").Value = 3555
$ws.Range("M126").Value = -1085

$ws.Range("H132").Value = 1366
$ws.Range("I132").Value = 1257.5
$ws.Range("J132").Value = 2288.25
$ws.Range("K132").Value = 3772.5
$ws.Range("L132").Value = 6864.75
$ws.Range("M132").Value = -1242.5
$ws.Range("N132").Value = -11924.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 221192.5
$ws.Range("I2").Value = 220100.2
$ws.Range("J2").Value = 222284.8
$ws.Range("K2").Value = 1320601.2
$ws.Range("L2").Value = 1333708.8
$ws.Range("M2").Value = -1320488.2
$ws.Range("N2").Value = -1333934.8

$ws.Range("H10").Value = 243.25
$ws.Range("I10").Value = 323.66666
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 970.9999799999999
$ws.Range("L10").Value = 6
$ws.Range("M10").Value = -831.9999799999999
$ws.Range("N10").Value = -284

$ws.Range("H137").Value = 2287.4119
$ws.Range("I137").Value = 1760.6666
$ws.Range("J137").Value = 2880
$ws.Range("K137").Value = 5281.9998
$ws.Range("L137").Value = 8640
$ws.Range("M137").Value = -181.9997999999996
$ws.Range("N137").Value = -18840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 281
$ws.Range("I2").Value = 309.44446
$ws.Range("K2").Value = 309.44446
$ws.Range("M2").Value = -196.44446

$ws.Range("H43").Value = 16993.4
$ws.Range("I43").Value = 11767
$ws.Range("J43").Value = 37899
$ws.Range("K43").Value = 11767
$ws.Range("L43").Value = 37899
$ws.Range("M43").Value = -11616
$ws.Range("N43").Value = -38201

$ws.Range("H80").Value = 999.6667
$ws.Range("J80").Value = 1100
$ws.Range("L80").Value = 1100
$ws.Range("N80").Value = -3096

$ws.Range("H83").Value = 999.6667
$ws.Range("J83").Value = 1100
$ws.Range("L83").Value = 5500
$ws.Range("N83").Value = -15484

$ws.Range("H97").Value = 233.27272
$ws.Range("I97").Value = 221.6
$ws.Range("K97").Value = 221.6
$ws.Range("M97").Value = 274.4

$ws.Range("H122").Value = 3972.4644
$ws.Range("I122").Value = 3180.0557
$ws.Range("J122").Value = 5398.8
$ws.Range("K122").Value = 9540.167099999999
$ws.Range("L122").Value = 16196.4
$ws.Range("M122").Value = -7090.167099999999
$ws.Range("N122").Value = -21096.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4220
$ws.Range("J40").Value = 4500
$ws.Range("L40").Value = 4500
$ws.Range("N40").Value = -4772

$ws.Range("H55").Value = 530.7143
$ws.Range("J55").Value = 699.4
$ws.Range("L55").Value = 699.4
$ws.Range("N55").Value = -1045.4

$ws.Range("H68").Value = 3001.5
$ws.Range("J68").Value = 3201.8
$ws.Range("L68").Value = 3201.8
$ws.Range("N68").Value = -4699.8

$ws.Range("H71").Value = 3001.5
$ws.Range("J71").Value = 3201.8
$ws.Range("L71").Value = 16009
$ws.Range("N71").Value = -23497

$ws.Range("H93").Value = 475.2
$ws.Range("I93").Value = 493
$ws.Range("J93").Value = 404
$ws.Range("K93").Value = 493
$ws.Range("L93").Value = 404
$ws.Range("M93").Value = 755
$ws.Range("N93").Value = -2900

$ws.Range("H99").Value = 100000
$ws.Range("J99").Value = 100000
$ws.Range("L99").Value = 100000
$ws.Range("N99").Value = -105990

$ws.Range("H122").Value = 7969.478
$ws.Range("I122").Value = 8353.467000000001
$ws.Range("K122").Value = 25060.401
$ws.Range("M122").Value = -22610.401

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 504.5
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1336
